# Verleiherrechnung Perfect Days 23.05.24 und Werbungskosten Schabi
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ausgaben")

# --- Prepare formatting for the new rows by copying from existing, similarly
# --- shaped rows, then overwrite the values (keeps style indices stable
# --- instead of minting new duplicate xf records). ---

# Row 45 ("Verleiher" entry) has the same column layout (incl. B/D dates) as row 44.
$ws.Range("A44:K44").Copy()
$ws.Range("A45:K45").PasteSpecial(-4122)

# Rows 46 & 47 ("Werbung" entries) have no B column - row 42 is a matching template.
$ws.Range("A42").Copy()
$ws.Range("A46").PasteSpecial(-4122)
$ws.Range("A47").PasteSpecial(-4122)
$ws.Range("C42:K42").Copy()
$ws.Range("C46:K46").PasteSpecial(-4122)
$ws.Range("C47:K47").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 45: DCM Film Distribution GmbH - Filmmiete "Perfect Days" ---
$ws.Range("A45").Value = "Verleiher"
$ws.Range("B45").Value = 45435
$ws.Range("C45").Value = "Film: Perfect Days"
$ws.Range("D45").Value = 45439
$ws.Range("E45").Value = 253.55
$ws.Range("F45").Value = "DCM Film Distribution GmbH"
$ws.Range("G45").Value = "Kreuzstrasse 2, 8008 Zürich"
$ws.Range("I45").Value = "03019739"
$ws.Range("H45").Value = "96 57660 00000 00000 00301 97392"
$ws.Range("J45").Value = "4404"
$ws.Range("K45").Value = "Filmmiete Kino"

# --- Row 46: Stefan Jablonski - Druckkosten Monatsprogramm ---
$ws.Range("A46").Value = "Werbung"
$ws.Range("C46").Value = "Druckkosten Monatsprogram Jan-Mai 24"
$ws.Range("D46").Value = 45435
$ws.Range("E46").Value = 479.5
$ws.Range("F46").Value = "Stefan Jablonski"
$ws.Range("G46").Value = "Sonnhaldenstrasse 13, 5734 Reinach AG"

# --- Row 47: Stefan Jablonski - Gestaltung Programm und Dias ---
$ws.Range("A47").Value = "Werbung"
$ws.Range("C47").Value = "Gestaltung Porgramm und Dias"
$ws.Range("D47").Value = 45435
$ws.Range("E47").Value = 420
$ws.Range("F47").Value = "Stefan Jablonski"
$ws.Range("G47").Value = "Sonnhaldenstrasse 13, 5734 Reinach AG"

$ws.Range("I46").Value = "001"
$ws.Range("I47").Value = "002"

$ws.Range("J46").Value = "4406"
$ws.Range("J46").Style = "Normal 2"
$ws.Range("J47").Value = "4406"
$ws.Range("J47").Style = "Normal 2"

$ws.Range("K46").Value = "Werbung Kino"
$ws.Range("K47").Value = "Werbung Kino"

# --- New small total box at K50:L50 (same "Normal 2" look as the other L/M boxes) ---
$ws.Range("K50:L50").Merge()
$ws.Range("K50:L50").Style = "Normal 2"

# --- Grow the table to include the three new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K47"))

# --- Match the author's final selection/cursor position ---
$ws.Activate()
$ws.Range("I47").Select()
